$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.479.88"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "1.688.66"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'316.24"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.3899"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").Value = "'0.4032"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'1.489"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'0.9997"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "'52.49"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "'0.08776"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'26.16"
$ws.Range("E13").Value = "  +11.46%  "
$ws.Range("D14").Value = "'7.496"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").Value = "'8.153"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "'0.00001349"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "1.686.53"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "'98.12"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "'0.07259"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").Value = "'19.97"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'7.282"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "'14.24"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").Value = "24.477.00"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").Value = "'3.047"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").Value = "'2.343"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "'22.61"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "'167.64"
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("D29").Value = "'8.640"
$ws.Range("E29").Value = "  +6.21%  "
$ws.Range("D30").Value = "'5.366"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("D31").Value = "'138.44"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "1.870.23"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "'0.08771"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'7.329"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'2.116"
$ws.Range("E35").Value = "  +6.92%  "
$ws.Range("D36").Value = "'1.048"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").Value = "'0.03025"
$ws.Range("E37").Value = "  +9.83%  "
$ws.Range("D38").Value = "'0.2791"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").Value = "'10.87"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D42").Value = "'14.19"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'1.473"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "'17.64"
$ws.Range("E44").Value = "  +11.86%  "
$ws.Range("D45").Value = "'2.644"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").Value = "'0.7256"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "'4.270"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'1.429"
$ws.Range("E48").Value = "  +8.80%  "
$ws.Range("D49").Value = "'0.9984"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'139.33"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").Value = "'0.08073"
$ws.Range("E51").Value = "  +0.88%  "

# Row 40/41: Stellar and TheSandbox swap places with refreshed price/volume data
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.8078"
$ws.Range("E40").Value = "  +4.82%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.09135"
$ws.Range("E41").Value = "  -0.90%  "
